$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.838.91'
$ws.Range('E2').Value = '  -0.74%  '
$ws.Range('D3').Value = '1.906.65'
$ws.Range('E3').Value = '  -0.17%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.34%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '313.38'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.94%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.003'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.27%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5032'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +4.09%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3816'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.11%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07278'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.17%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9081'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.89%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.84'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.07%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07676'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.23%  '
$ws.Range('D13').Value = '1.908.40'
$ws.Range('E13').Value = '  +0.10%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.478'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.47%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '91.74'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.16%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.003'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.31%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008714'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.53%  '
$ws.Range('E18').Value = '  -0.24%  '
$ws.Range('D19').Value = '27.869.43'
$ws.Range('E19').Value = '  -0.71%  '
$ws.Range('E20').Value = '  -1.64%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.169'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.34%  '
$ws.Range('E22').Value = '  -0.90%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.584'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.71%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '154.13'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.21%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.881'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.19%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.212'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +4.01%  '
$ws.Range('E27').Value = '  -0.78%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '115.39'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.35%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.899'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.51%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.09004'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.51%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.210'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.20%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.224'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.40%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.656'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.53%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7613'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.85%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.02062'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.04%  '
$ws.Range('B36').Value = 'TrustWalletToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.095'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.48%  '
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.484'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -5.86%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5524'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.62%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.012'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.62%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.05246'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.22%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.878'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.27%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.470'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.61%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1510'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.35%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '110.89'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.87%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.58'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.62%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4812'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.49%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.003'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.24%  '
$ws.Range('E48').Value = '  -1.70%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '67.29'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.07%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06063'
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.9024'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.21%  '
